# Trade #18 closed at 2026-02-17 04:07:59 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly closed trade #18 for the MarketMaking strategy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: bump running totals now that trade #18 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.08   # Current Capital
$summary.Range("B4").Value = 0.08      # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 18        # Total Trades
$summary.Range("B7").Value = 6         # Winning Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: same rollup, but scoped to the MarketMaking row.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.08     # Capital
$status.Range("D4").Value = 18         # Trades
$status.Range("E4").Value = 0.08       # P&L $
$status.Range("F4").Value = 0.08       # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------------
# Append the new trade row (#18) to both "All Trades" and "MarketMaking".
# ---------------------------------------------------------------------------
$tradeRow = @{
    A = 18
    B = "2026-02-17"
    C = "04:07:54"
    D = "MarketMaking"
    E = "UP"
    F = 0.19
    G = 0.26
    H = "CLOSED"
    I = 36.8421
    J = 0.07000000000000001
    K = 100.08
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.12
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A19").Value = $tradeRow.A

    # The Date column holds plain text like "2026-02-17" (matching every
    # other row in the column), not a real Excel date serial. Force the
    # cell to text first so the ISO-looking string isn't auto-converted,
    # then restore the default "Normal" style so no formatting lingers on
    # the cell itself.
    $dateCell = $ws.Range("B19")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $tradeRow.B
    $dateCell.Style = "Normal"

    $ws.Range("C19").Value = $tradeRow.C
    $ws.Range("D19").Value = $tradeRow.D
    $ws.Range("E19").Value = $tradeRow.E
    $ws.Range("F19").Value = $tradeRow.F
    $ws.Range("G19").Value = $tradeRow.G
    $ws.Range("H19").Value = $tradeRow.H
    $ws.Range("I19").Value = $tradeRow.I
    $ws.Range("J19").Value = $tradeRow.J
    $ws.Range("K19").Value = $tradeRow.K
    $ws.Range("L19").Value = $tradeRow.L
    $ws.Range("M19").Value = $tradeRow.M
    $ws.Range("N19").Value = $tradeRow.N
    $ws.Range("O19").Value = $tradeRow.O
    $ws.Range("P19").Value = $tradeRow.P
    $ws.Range("Q19").Value = $tradeRow.Q
}
